# Auto-generated edit script applying the diff between before.xlsx and target cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.856.56"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.780.80"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.58%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -5.50%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "2.037.92"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.828.55"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "33.855.77"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").Value = "1.390.69"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.37%  "
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0509"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("E47").Value = "  +7.90%  "
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "1.939.05"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("E51").Value = "  +0.04%  "
